$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh.
# For "Price" (column D) cells whose new value parses as a plain number,
# force the cell to text format first so Excel keeps it as a literal
# string (e.g. "584.45") instead of silently coercing it to a Double -
# matching the text-valued cells already used throughout this sheet.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '66.704.37'
$ws.Range("E2").Value = '  -1.12%  '

$ws.Range("D3").Value = '3.517.53'
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '584.45'
$ws.Range("E5").Value = '  -2.47%  '

$ws.Range("D6").Value = '175.59'
$ws.Range("E6").Value = '  -2.58%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '3.514.67'
$ws.Range("E8").Value = '  +0.46%  '

$ws.Range("D9").Value = '0.597'
$ws.Range("E9").Value = '  -2.17%  '

$ws.Range("E10").Value = '  -2.68%  '

$ws.Range("E11").Value = '  -1.99%  '

$ws.Range("D12").Value = '0.424'
$ws.Range("E12").Value = '  -3.02%  '

$ws.Range("D13").Value = '4.119.84'
$ws.Range("E13").Value = '  +0.50%  '

$ws.Range("E14").Value = '  -5.79%  '

$ws.Range("E15").Value = '  -1.61%  '

$ws.Range("D16").Value = '66.698.00'
$ws.Range("E16").Value = '  -1.08%  '

$ws.Range("E17").Value = '  -2.49%  '

$ws.Range("D18").Value = '3.504.98'
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("E19").Value = '  -3.78%  '

$ws.Range("D20").Value = '14.01'
$ws.Range("E20").Value = '  -2.00%  '

$ws.Range("D21").Value = '382.15'
$ws.Range("E21").Value = '  -2.18%  '

$ws.Range("D22").Value = '7.91'
$ws.Range("E22").Value = '  -0.57%  '

$ws.Range("D23").Value = '0.550'
$ws.Range("E23").Value = '  +1.50%  '

$ws.Range("E24").Value = '  +0.23%  '

$ws.Range("D25").Value = '72.37'
$ws.Range("E25").Value = '  -2.15%  '

$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("E27").Value = '  -1.24%  '

$ws.Range("D28").Value = '9.89'
$ws.Range("E28").Value = '  -4.76%  '

$ws.Range("E29").Value = '  -1.70%  '

$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("D31").Value = '24.61'
$ws.Range("E31").Value = '  +4.42%  '

$ws.Range("E32").Value = '  -4.55%  '

$ws.Range("E33").Value = '  -2.78%  '

$ws.Range("E34").Value = '  -5.58%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '7.27'
$ws.Range("E35").Value = '  -1.77%  '

$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("E37").Value = '  -1.78%  '

$ws.Range("D38").Value = '30.12'
$ws.Range("E38").Value = '  +13.67%  '

$ws.Range("D39").Value = '161.61'
$ws.Range("E39").Value = '  -0.90%  '

$ws.Range("D40").Value = '0.899'
$ws.Range("E40").Value = '  +3.20%  '

$ws.Range("E41").Value = '  -5.01%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '4.51'
$ws.Range("E42").Value = '  -3.03%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '6.52'
$ws.Range("E43").Value = '  -4.85%  '

$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '2.55'
$ws.Range("E44").Value = '  -9.56%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.729.86'
$ws.Range("E45").Value = '  -4.50%  '

$ws.Range("E46").Value = '  -2.94%  '

$ws.Range("D47").Value = '40.75'
$ws.Range("E47").Value = '  -2.30%  '

$ws.Range("D48").Value = '25.07'
$ws.Range("E48").Value = '  -6.60%  '

$ws.Range("E49").Value = '  -2.61%  '

$ws.Range("D50").Value = '325.29'
$ws.Range("E50").Value = '  -2.48%  '

$ws.Range("E51").Value = '  -3.70%  '
